$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Programa resumido (PT)
$ok1 = $d.Content.Find.Execute('1.Destilação 2.Absorção3.Extração líquido-líquido4.Adsorção5.Cristalização', $false, $false, $false, $false, $false, $true, 1, $false, '1.Destilação ^l2.Absorção^l3.Extração líquido-líquido^l4.Adsorção^l5.Cristalização', 2)
Write-Output "replace1: $ok1"

# Programa (PT)
$ok2 = $d.Content.Find.Execute('1) Equilíbrio líquido-vapor (Tempo estimado: 2 horas); Separação simples: Destilação flash (Tempo estimado: 2 horas e Destilação Diferencial (Tempo estimado: 2 horas); Destilação contínua (Retificação): Método de McCabe-Thiele; Eficiência de estágio e eficiência global (Tempo estimado: 16 horas); Destilação multicomponentes – método FUG (Tempo estimado: 8 horas);2) Absorção e dessorção: tipos de torres; Solubilidade de gases em líquidos; Operações em paralelo e contracorrente; Taxas de transferência de massa; Operações multiestágios em contracorrente (Tempo estimado: 8 horas);3) Extração líquido-líquido: equilíbrio líquido-líquido; Extração em estágio único e em múltiplos estágios; Coeficientes de distribuição (Tempo estimado: 14 horas);4) Adsorção: fundamentos; Operações em único estágio e em contato contínuo  (Tempo estimado: 4 horas); 5) Cristalização: Caracterização de partículas e Projeto de Cristalizadores (Tempo estimado: 4 horas).', $false, $false, $false, $false, $false, $true, 1, $false, '1) Equilíbrio líquido-vapor (Tempo estimado: 2 horas); Separação simples: Destilação flash (Tempo estimado: 2 horas e Destilação Diferencial (Tempo estimado: 2 horas); Destilação contínua (Retificação): Método de McCabe-Thiele; Eficiência de estágio e eficiência global (Tempo estimado: 16 horas); Destilação multicomponentes – método FUG (Tempo estimado: 8 horas);^l2) Absorção e dessorção: tipos de torres; Solubilidade de gases em líquidos; Operações em paralelo e contracorrente; Taxas de transferência de massa; Operações multiestágios em contracorrente (Tempo estimado: 8 horas);^l3) Extração líquido-líquido: equilíbrio líquido-líquido; Extração em estágio único e em múltiplos estágios; Coeficientes de distribuição (Tempo estimado: 14 horas);^l4) Adsorção: fundamentos; Operações em único estágio e em contato contínuo  (Tempo estimado: 4 horas); ^l5) Cristalização: Caracterização de partículas e Projeto de Cristalizadores (Tempo estimado: 4 horas).', 2)
Write-Output "replace2: $ok2"

# Programa (EN)
$ok3 = $d.Content.Find.Execute('1) Vapor-liquid equilibrium (Estimated time: 2 hours); Simple separation: Flash distillation (Estimated time: 2 hours and Differential distillation (Estimated time: 2 hours); Continuous distillation (Rectification): McCabe-Thiele method; Stage efficiency and overall efficiency (Estimated time: 16 hours); Multicomponent distillation – FUG method (Estimated time: 8 hours);2) Absorption and desorption: types of towers; Solubility of gases in liquids; Parallel and countercurrent operations; Mass transfer rates; Countercurrent multistage operations (Estimated time: 8 hours);3) Liquid-liquid extraction: liquid-liquid balance; Single-stage and multi-stage extraction; Distribution coefficients (Estimated time: 14 hours);4) Adsorption: fundamentals; Single stage and continuous contact operations (Estimated time: 4 hours);5) Crystallization: Characterization of particles and Design of Crystallizers (Estimated time: 4 hours).', $false, $false, $false, $false, $false, $true, 1, $false, '1) Vapor-liquid equilibrium (Estimated time: 2 hours); Simple separation: Flash distillation (Estimated time: 2 hours and Differential distillation (Estimated time: 2 hours); Continuous distillation (Rectification): McCabe-Thiele method; Stage efficiency and overall efficiency (Estimated time: 16 hours); Multicomponent distillation – FUG method (Estimated time: 8 hours);^l2) Absorption and desorption: types of towers; Solubility of gases in liquids; Parallel and countercurrent operations; Mass transfer rates; Countercurrent multistage operations (Estimated time: 8 hours);^l3) Liquid-liquid extraction: liquid-liquid balance; Single-stage and multi-stage extraction; Distribution coefficients (Estimated time: 14 hours);^l4) Adsorption: fundamentals; Single stage and continuous contact operations (Estimated time: 4 hours);^l5) Crystallization: Characterization of particles and Design of Crystallizers (Estimated time: 4 hours).', 2)
Write-Output "replace3: $ok3"

# Bibliografia
$ok4 = $d.Content.Find.Execute('1) TREYBAL, R. E. Mass-Transfer Operations. 3ed. Auckland: McGraw-Hill, 784p. 1980;2)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;3)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005.Bibliografia Complementar:1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v. 2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;2)COULSON & Richardson''s Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;4)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;5)PERRY''s chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008;6)SEADER, J. D; HENLEY, E. J. Separation Process Principles. 2ed. Hoboken, N.J: Wiley, 756p. 2006.', $false, $false, $false, $false, $false, $true, 1, $false, '1) TREYBAL, R. E. Mass-Transfer Operations. 3ed. Auckland: McGraw-Hill, 784p. 1980;^l2)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;^l3)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005.^l^lBibliografia Complementar:^l1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v. 2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;^l2)COULSON & Richardson''s Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;^l3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;^l4)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;^l5)PERRY''s chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008;^l6)SEADER, J. D; HENLEY, E. J. Separation Process Principles. 2ed. Hoboken, N.J: Wiley, 756p. 2006.', 2)
Write-Output "replace4: $ok4"
